$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (74) down into the
# new rows 75-83 so the new rows inherit the same fonts / wrap / alignment
# styles used throughout the table (style 1 for column A, style 2 for B:D).
$ws.Range("A74:D74").Copy($ws.Range("A75:D83")) | Out-Null

# New vocabulary rows: word, definition, example 1, example 2
$rows = @(
    @{ Row = 75; Height = 60; A = "beam";        B = "a long heavy piece of wood or metal used in building houses, bridges etc"; C = "Workers used steal beams to brace the roof."; D = "Her head hit the beam and she slid down and hit her shoulder." }
    @{ Row = 76; Height = 60; A = "brace";       B = "brace yourself, to make something stronger by supporting it"; C = "The carpenter uses a brace to hold pieces of wood in place"; D = "Nancy braced herself for the inevitable arguments." }
    @{ Row = 77; Height = 75; A = "classy";      B = "stylish and sophisticated."; C = "She took us to a very classy seafood restaurant in the old part of the city"; D = "The hotel is classy but relaxed." }
    @{ Row = 78; Height = 45; A = "consolation"; B = "the comfort received by a person after a loss or disappointment."; C = "your company has been a great consolation to me."; D = "The Church was the main consolation in a short and hard life." }
    @{ Row = 79; Height = 75; A = "successor";   B = "someone who takes a job or position previously held by someone else"; C = "someone who takes a job or position previously held by someone else"; D = "His successor died after only 15 months in office." }
    @{ Row = 80; Height = 60; A = "subtle";      B = "not easy to notice or understand unless you pay careful attention"; C = " you have to be some more subtle with people."; D = "The pictures are similar, but there are subtle differences between them." }
    @{ Row = 81; Height = 60; A = "stubborn";    B = "determined not to change your mind, even when people think you are being unreasonable"; C = "he was too stubborn to admit that he was worng"; D = "Why are you so stubborn?" }
    @{ Row = 82; Height = 60; A = "admit";       B = "to agree unwillingly that something is true or that someone else is right"; C = "You may not like her, but you have to admit that she’s good at her job."; D = "I must admit, I didn’t actually do anything to help her." }
    @{ Row = 83; Height = 45; A = "spare";       B = "not being used or not needed at the present time"; C = "we have a spare in the trunk."; D = "a spare bass guitar line." }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Rows($r.Row).RowHeight = $r.Height
}

# Update the view: scroll so the newly added rows are visible and move the
# active selection down near the new data, mirroring the author's edit.
$win = $excel.ActiveWindow
$win.ScrollRow = 78
$win.ScrollColumn = 1
$ws.Range("B89").Select() | Out-Null
